$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the sample credentials row: the "Usuario" (username) value now matches
# the e-mail address, and a new password is used for both Senha/ConfirmarSenha.
$ws.Range("A2").Value = "saulo.silva@rsinet.com.br"
$ws.Range("C2").Value = "manodoCeu12"
$ws.Range("D2").Value = "manodoCeu12"

# Columns got wider content, so resize them to fit.
$ws.Columns("A").AutoFit() | Out-Null
$ws.Columns("C").AutoFit() | Out-Null

$wb.Save()
